# GUI: Updated the statistics.
#
# The four numbers below are the "raw" test-run totals that feed this
# summary sheet; every other changed cell (E2, I2, L2, L3, L5, L6, L7,
# N2, N3, P3 ...) is a formula that recomputes automatically once these
# inputs change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 1     # was 2
$ws.Range("C2").Value2 = 24    # was 23
$ws.Range("G2").Value2 = 193   # was 176
$ws.Range("H2").Value2 = 250   # was 243

$wb.Application.CalculateFull()
